$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" '66.354.73'
Set-TextValue "E2" '  +1.50%  '
Set-TextValue "D3" '3.619.08'
Set-TextValue "E3" '  +2.04%  '
Set-TextValue "D4" '0.998'
Set-TextValue "E4" '  -0.23%  '
Set-TextValue "D5" '606.97'
Set-TextValue "E5" '  +0.99%  '
Set-TextValue "D6" '140.12'
Set-TextValue "E6" '  +1.23%  '
Set-TextValue "D7" '3.618.88'
Set-TextValue "E7" '  +2.11%  '
Set-TextValue "D8" '1.00'
Set-TextValue "E8" '  -0.04%  '
Set-TextValue "D9" '0.500'
Set-TextValue "E9" '  +1.44%  '
Set-TextValue "D10" '0.128'
Set-TextValue "E10" '  +2.39%  '
Set-TextValue "D11" '7.28'
Set-TextValue "E11" '  +6.14%  '
Set-TextValue "D12" '0.396'
Set-TextValue "E12" '  +2.36%  '
Set-TextValue "D13" '4.218.85'
Set-TextValue "E13" '  +1.60%  '
Set-TextValue "D14" '28.81'
Set-TextValue "E14" '  +5.59%  '
Set-TextValue "D15" '0.0000189'
Set-TextValue "E15" '  +2.51%  '
Set-TextValue "D16" '3.605.35'
Set-TextValue "E16" '  +1.43%  '
Set-TextValue "E17" '  +0.02%  '
Set-TextValue "D18" '66.366.22'
Set-TextValue "E18" '  +1.64%  '
Set-TextValue "D19" '10.22'
Set-TextValue "E19" '  -0.39%  '
Set-TextValue "D20" '14.81'
Set-TextValue "E20" '  +3.40%  '
Set-TextValue "D21" '5.95'
Set-TextValue "E21" '  -0.09%  '
Set-TextValue "D22" '400.88'
Set-TextValue "E22" '  +1.84%  '
Set-TextValue "D23" '0.596'
Set-TextValue "E23" '  +3.60%  '
Set-TextValue "D24" '3.763.59'
Set-TextValue "E24" '  +1.87%  '
Set-TextValue "D25" '74.77'
Set-TextValue "E25" '  +1.26%  '
Set-TextValue "D26" '0.999'
Set-TextValue "E26" '  +0.03%  '
Set-TextValue "D27" '0.0000121'
Set-TextValue "E27" '  +4.10%  '
Set-TextValue "D28" '8.25'
Set-TextValue "E28" '  +5.75%  '
Set-TextValue "D29" '1.66'
Set-TextValue "E29" '  +28.92%  '
Set-TextValue "B30" 'InternetComputer(DFINITY)'
Set-TextValue "C30" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D30" '8.73'
Set-TextValue "E30" '  +7.02%  '
Set-TextValue "B31" 'PancakeSwap'
Set-TextValue "C31" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D31" '2.36'
Set-TextValue "E31" '  +3.10%  '
Set-TextValue "E32" '  -0.09%  '
Set-TextValue "D33" '3.619.99'
Set-TextValue "E33" '  +1.61%  '
Set-TextValue "D34" '24.78'
Set-TextValue "E34" '  +3.85%  '
Set-TextValue "D35" '0.150'
Set-TextValue "E35" '  +2.90%  '
Set-TextValue "E36" '  +0.02%  '
Set-TextValue "D37" '5.48'
Set-TextValue "E37" '  +9.30%  '
Set-TextValue "D38" '1.66'
Set-TextValue "E38" '  +6.26%  '
Set-TextValue "D39" '7.15'
Set-TextValue "E39" '  +2.47%  '
Set-TextValue "D40" '169.31'
Set-TextValue "E40" '  -0.04%  '
Set-TextValue "D41" '0.0848'
Set-TextValue "E41" '  +5.45%  '
Set-TextValue "D42" '0.849'
Set-TextValue "E42" '  +2.65%  '
Set-TextValue "D43" '27.17'
Set-TextValue "E43" '  +2.83%  '
Set-TextValue "D44" '1.30'
Set-TextValue "E44" '  +8.16%  '
Set-TextValue "B45" 'Filecoin'
Set-TextValue "C45" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D45" '4.60'
Set-TextValue "E45" '  +3.79%  '
Set-TextValue "B46" 'OKB'
Set-TextValue "C46" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D46" '43.12'
Set-TextValue "E46" '  +0.86%  '
Set-TextValue "B47" 'Stacks'
Set-TextValue "C47" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D47" '1.74'
Set-TextValue "E47" '  +3.63%  '
Set-TextValue "B48" 'FirstDigitalUSD'
Set-TextValue "C48" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D48" '0.998'
Set-TextValue "E48" '  -0.35%  '
Set-TextValue "D49" '7.08'
Set-TextValue "E49" '  +4.02%  '
Set-TextValue "D50" '2.482.02'
Set-TextValue "E50" '  +3.08%  '
Set-TextValue "D51" '0.919'
Set-TextValue "E51" '  +10.43%  '
